$d = $word.ActiveDocument

# --- Remove the existing "_GoBack" bookmark that currently sits between
#     " on" and " trees" in the "Table XX." caption paragraph. ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Place a temporary barrier bookmark right after the run that ends
#     "...simulated under time variable speciation rates" so that the
#     text edit below cannot coalesce that run with the one that follows
#     it (bookmarks act as hard boundaries against run merging). ---
$boundary = $d.Content
$boundary.Find.Execute(" simulated under time variable speciation rates", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$boundaryPoint = $d.Range($boundary.End, $boundary.End)
$d.Bookmarks.Add("TempBarrier", $boundaryPoint) | Out-Null

# --- Delete the sentence "Tree type refers to the extinction fraction
#     for the birth-death trees. " leaving ". The rejection rate..." ---
$toDelete = $d.Content
$toDelete.Find.Execute("Tree type refers to the extinction fraction for the birth-death trees. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$toDelete.Text = ""

# --- Re-insert the "_GoBack" bookmark immediately before "The rejection
#     rate", matching the cursor position left by the author's edit. ---
$target = $d.Content
$target.Find.Execute("The rejection rate is the proportion", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPoint = $d.Range($target.Start, $target.Start)
$d.Bookmarks.Add("_GoBack", $insertPoint) | Out-Null

# --- Remove the temporary barrier bookmark now that the runs are split
#     the way we need them. ---
$d.Bookmarks("TempBarrier").Delete()
